$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.384145666666667
$ws.Range("H2").Value = 4.152437
$ws.Range("I2").Value = 0.1014617184198512
$ws.Range("J2").Value = 0.1334061399754118
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 109.1447706666667
$ws.Range("N2").Value = 327.434312
$ws.Range("O2").Value = 0.3535542089399963
$ws.Range("P2").Value = 0.3655959674582361
$ws.Range("Q2").Value = 151.0722613575938
$ws.Range("R2").Value = 1359.650352218344
$ws.Range("S2").Value = 0.03587221759362313
$ws.Range("T2").Value = 0.04877274680917953

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.384145666666667
$ws.Range("H3").Value = 4.152437
$ws.Range("I3").Value = 0.1014617184198512
$ws.Range("J3").Value = 0.1334061399754118
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.980825
$ws.Range("N3").Value = 143.942475
$ws.Range("O3").Value = 0.155424969272891
$ws.Range("P3").Value = 0.1607186127944892
$ws.Range("Q3").Value = 66.41245100684168
$ws.Range("R3").Value = 597.712059061575
$ws.Range("S3").Value = 0.01576968446778008
$ws.Range("T3").Value = 0.02144084975511563

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.384145666666667
$ws.Range("H4").Value = 4.152437
$ws.Range("I4").Value = 0.1014617184198512
$ws.Range("J4").Value = 0.1334061399754118
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 62.26741999999999
$ws.Range("N4").Value = 186.80226
$ws.Range("O4").Value = 0.2017037397794264
$ws.Range("P4").Value = 0.2085735992386923
$ws.Range("Q4").Value = 86.18717956751333
$ws.Range("R4").Value = 775.6846161076199
$ws.Range("S4").Value = 0.0204652080497311
$ws.Range("T4").Value = 0.02782499877521243

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.384145666666667
$ws.Range("H5").Value = 4.152437
$ws.Range("I5").Value = 0.1014617184198512
$ws.Range("J5").Value = 0.1334061399754118
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 58.81030666666667
$ws.Range("N5").Value = 176.43092
$ws.Range("O5").Value = 0.1905050633580386
$ws.Range("P5").Value = 0.1969935053322898
$ws.Range("Q5").Value = 81.40203112800445
$ws.Range("R5").Value = 732.61828015204
$ws.Range("S5").Value = 0.01932897109598922
$ws.Range("T5").Value = 0.02628014314660649

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.384145666666667
$ws.Range("H6").Value = 4.152437
$ws.Range("I6").Value = 0.1014617184198512
$ws.Range("J6").Value = 0.1334061399754118
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 30.503993
$ws.Range("N6").Value = 61.007986
$ws.Range("O6").Value = 0.09881201864964768
$ws.Range("P6").Value = 0.06811831517629259
$ws.Range("Q6").Value = 42.22196972698033
$ws.Range("R6").Value = 253.331818361882
$ws.Range("S6").Value = 0.01002563721272764
$ws.Range("T6").Value = 0.009087401489297707

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.458038666666667
$ws.Range("H7").Value = 7.374116000000001
$ws.Range("I7").Value = 0.1801810554109116
$ws.Range("J7").Value = 0.2369096391566985
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 109.1447706666667
$ws.Range("N7").Value = 327.434312
$ws.Range("O7").Value = 0.3535542089399963
$ws.Range("P7").Value = 0.3655959674582361
$ws.Range("Q7").Value = 268.2820665631324
$ws.Range("R7").Value = 2414.538599068192
$ws.Range("S7").Value = 0.06370377051177849
$ws.Range("T7").Value = 0.08661320872767481

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.458038666666667
$ws.Range("H8").Value = 7.374116000000001
$ws.Range("I8").Value = 0.1801810554109116
$ws.Range("J8").Value = 0.2369096391566985
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 47.980825
$ws.Range("N8").Value = 143.942475
$ws.Range("O8").Value = 0.155424969272891
$ws.Range("P8").Value = 0.1607186127944892
$ws.Range("Q8").Value = 117.9387231085667
$ws.Range("R8").Value = 1061.4485079771
$ws.Range("S8").Value = 0.028004635000798
$ws.Range("T8").Value = 0.03807578856290758

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.458038666666667
$ws.Range("H9").Value = 7.374116000000001
$ws.Range("I9").Value = 0.1801810554109116
$ws.Range("J9").Value = 0.2369096391566985
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 62.26741999999999
$ws.Range("N9").Value = 186.80226
$ws.Range("O9").Value = 0.2017037397794264
$ws.Range("P9").Value = 0.2085735992386923
$ws.Range("Q9").Value = 153.0557260335733
$ws.Range("R9").Value = 1377.50153430216
$ws.Range("S9").Value = 0.03634319271378492
$ws.Range("T9").Value = 0.04941309613325245

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.458038666666667
$ws.Range("H10").Value = 7.374116000000001
$ws.Range("I10").Value = 0.1801810554109116
$ws.Range("J10").Value = 0.2369096391566985
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 58.81030666666667
$ws.Range("N10").Value = 176.43092
$ws.Range("O10").Value = 0.1905050633580386
$ws.Range("P10").Value = 0.1969935053322898
$ws.Range("Q10").Value = 144.5580077851911
$ws.Range("R10").Value = 1301.02207006672
$ws.Range("S10").Value = 0.03432540337697397
$ws.Range("T10").Value = 0.04666966026448596

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.458038666666667
$ws.Range("H11").Value = 7.374116000000001
$ws.Range("I11").Value = 0.1801810554109116
$ws.Range("J11").Value = 0.2369096391566985
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 30.503993
$ws.Range("N11").Value = 61.007986
$ws.Range("O11").Value = 0.09881201864964768
$ws.Range("P11").Value = 0.06811831517629259
$ws.Range("Q11").Value = 74.97999428172935
$ws.Range("R11").Value = 449.8799656903761
$ws.Range("S11").Value = 0.0178040538075762
$ws.Range("T11").Value = 0.01613788546837774

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 9.799863999999999
$ws.Range("H12").Value = 19.599728
$ws.Range("I12").Value = 0.7183572261692373
$ws.Range("J12").Value = 0.6296842208678898
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 109.1447706666667
$ws.Range("N12").Value = 327.434312
$ws.Range("O12").Value = 0.3535542089399963
$ws.Range("P12").Value = 0.3655959674582361
$ws.Range("Q12").Value = 1069.603908844522
$ws.Range("R12").Value = 6417.623453067135
$ws.Range("S12").Value = 0.2539782208345947
$ws.Range("T12").Value = 0.2302100119213818

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 9.799863999999999
$ws.Range("H13").Value = 19.599728
$ws.Range("I13").Value = 0.7183572261692373
$ws.Range("J13").Value = 0.6296842208678898
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 47.980825
$ws.Range("N13").Value = 143.942475
$ws.Range("O13").Value = 0.155424969272891
$ws.Range("P13").Value = 0.1607186127944892
$ws.Range("Q13").Value = 470.2055596078
$ws.Range("R13").Value = 2821.2333576468
$ws.Range("S13").Value = 0.1116506498043129
$ws.Range("T13").Value = 0.101201974476466

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 9.799863999999999
$ws.Range("H14").Value = 19.599728
$ws.Range("I14").Value = 0.7183572261692373
$ws.Range("J14").Value = 0.6296842208678898
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 62.26741999999999
$ws.Range("N14").Value = 186.80226
$ws.Range("O14").Value = 0.2017037397794264
$ws.Range("P14").Value = 0.2085735992386923
$ws.Range("Q14").Value = 610.2122476308799
$ws.Range("R14").Value = 3661.27348578528
$ws.Range("S14").Value = 0.1448953390159104
$ws.Range("T14").Value = 0.1313355043302275

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 9.799863999999999
$ws.Range("H15").Value = 19.599728
$ws.Range("I15").Value = 0.7183572261692373
$ws.Range("J15").Value = 0.6296842208678898
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 58.81030666666667
$ws.Range("N15").Value = 176.43092
$ws.Range("O15").Value = 0.1905050633580386
$ws.Range("P15").Value = 0.1969935053322898
$ws.Range("Q15").Value = 576.3330071316267
$ws.Range("R15").Value = 3457.99804278976
$ws.Range("S15").Value = 0.1368506888850754
$ws.Range("T15").Value = 0.1240437019211974

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 9.799863999999999
$ws.Range("H16").Value = 19.599728
$ws.Range("I16").Value = 0.7183572261692373
$ws.Range("J16").Value = 0.6296842208678898
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 30.503993
$ws.Range("N16").Value = 61.007986
$ws.Range("O16").Value = 0.09881201864964768
$ws.Range("P16").Value = 0.06811831517629259
$ws.Range("Q16").Value = 298.934982856952
$ws.Range("R16").Value = 1195.739931427808
$ws.Range("S16").Value = 0.07098232762934385
$ws.Range("T16").Value = 0.07455013442405678
